# Weekly update: insert a new daily record as a new row before the
# existing row 174, pushing the rest of the data set down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 174; this shifts the existing rows
# 174..242 down to 175..243 and extends the used range to row 243.
$ws.Rows.Item(174).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(174, 1).Value  = 7
$ws.Cells.Item(174, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(174, 3).Value  = "Ñuble"
$ws.Cells.Item(174, 4).Value  = 44755
$ws.Cells.Item(174, 5).Value  = 16
$ws.Cells.Item(174, 6).Value  = 100112009
$ws.Cells.Item(174, 7).Value  = "Acelga"
$ws.Cells.Item(174, 8).Value  = "Sin especificar"
$ws.Cells.Item(174, 9).Value  = "Segunda"
$ws.Cells.Item(174, 10).Value = 120
$ws.Cells.Item(174, 11).Value = 500
$ws.Cells.Item(174, 12).Value = 500
$ws.Cells.Item(174, 13).Value = 500
$ws.Cells.Item(174, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(174, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(174, 16).Value = 500
$ws.Cells.Item(174, 17).Value = 1
$ws.Cells.Item(174, 18).Value = "Hortaliza"
